$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.630.71"
$ws.Range("E2").Value = "  -1.41%  "
$ws.Range("D3").Value = "1.593.88"
$ws.Range("E3").Value = "  -1.68%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.01"
$ws.Range("E5").Value = "  -1.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.513"
$ws.Range("E6").Value = "  -1.06%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E9").Value = "  -1.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.63"
$ws.Range("E10").Value = "  -2.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0835"
$ws.Range("E11").Value = "  -1.42%  "
$ws.Range("D12").Value = "1.818.88"
$ws.Range("E12").Value = "  -1.59%  "
$ws.Range("D13").Value = "1.598.23"
$ws.Range("E13").Value = "  -1.21%  "
$ws.Range("E14").Value = "  -2.72%  "
$ws.Range("E15").Value = "  -2.90%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.13"
$ws.Range("E16").Value = "  +0.45%  "
$ws.Range("D17").Value = "26.611.35"
$ws.Range("E17").Value = "  -1.39%  "
$ws.Range("D18").Value = "0.0₃0730"
$ws.Range("E18").Value = "  -2.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "209.82"
$ws.Range("E19").Value = "  -1.82%  "
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.69"
$ws.Range("E21").Value = "  -2.20%  "
$ws.Range("E23").Value = "  -2.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.87"
$ws.Range("E24").Value = "  -2.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.40"
$ws.Range("E25").Value = "  -1.23%  "
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.14"
$ws.Range("E27").Value = "  -3.09%  "
$ws.Range("E28").Value = "  -1.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.32"
$ws.Range("E29").Value = "  -1.50%  "
$ws.Range("E30").Value = "  -1.63%  "
$ws.Range("E31").Value = "  -1.48%  "
$ws.Range("E32").Value = "  -3.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.665"
$ws.Range("E33").Value = "  -13.96%  "
$ws.Range("D35").Value = "1.293.25"
$ws.Range("E35").Value = "  -4.15%  "
$ws.Range("E36").Value = "  -0.42%  "
$ws.Range("E37").Value = "  -4.94%  "
$ws.Range("E38").Value = "  -3.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.833"
$ws.Range("E39").Value = "  -1.48%  "
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("E41").Value = "  -1.06%  "
$ws.Range("E42").Value = "  +0.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.19"
$ws.Range("E43").Value = "  -1.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.56"
$ws.Range("E44").Value = "  -2.23%  "
$ws.Range("D45").Value = "1.730.55"
$ws.Range("E45").Value = "  -1.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "89.45"
$ws.Range("E46").Value = "  -0.60%  "
$ws.Range("E47").Value = "  -1.89%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.830"
$ws.Range("E48").Value = "  -5.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0987"
$ws.Range("E49").Value = "  -3.35%  "
$ws.Range("E50").Value = "  -2.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.52"
$ws.Range("E51").Value = "  -1.77%  "
